$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the loopback length input value; dependent formulas (C2:G2, C3:G3) recalc automatically.
$ws.Range("A2").Value = 16000

# Move the active cell selection to F2, matching the saved view state.
$ws.Activate()
$ws.Range("F2").Select()
